$d = $word.ActiveDocument

$d.Content.Find.Execute("742×3=2226", $true, $false, $false, $false, $false, $true, 1, $false, "237×7=1659", 2) | Out-Null
$d.Content.Find.Execute("246×4=984", $true, $false, $false, $false, $false, $true, 1, $false, "621×5=3105", 2) | Out-Null
$d.Content.Find.Execute("844×3=2532", $true, $false, $false, $false, $false, $true, 1, $false, "659×6=3954", 2) | Out-Null
$d.Content.Find.Execute("724×8=5792", $true, $false, $false, $false, $false, $true, 1, $false, "523×4=2092", 2) | Out-Null
$d.Content.Find.Execute("147×3=441", $true, $false, $false, $false, $false, $true, 1, $false, "206×9=1854", 2) | Out-Null
$d.Content.Find.Execute("567×9=5103", $true, $false, $false, $false, $false, $true, 1, $false, "795×4=3180", 2) | Out-Null
$d.Content.Find.Execute("879×3=2637", $true, $false, $false, $false, $false, $true, 1, $false, "773×5=3865", 2) | Out-Null
$d.Content.Find.Execute("999×8=7992", $true, $false, $false, $false, $false, $true, 1, $false, "969×2=1938", 2) | Out-Null
$d.Content.Find.Execute("742×5=3710", $true, $false, $false, $false, $false, $true, 1, $false, "198×5=990", 2) | Out-Null
$d.Content.Find.Execute("571×3=1713", $true, $false, $false, $false, $false, $true, 1, $false, "646×5=3230", 2) | Out-Null
$d.Content.Find.Execute("312×5=1560", $true, $false, $false, $false, $false, $true, 1, $false, "219×3=657", 2) | Out-Null
$d.Content.Find.Execute("457×3=1371", $true, $false, $false, $false, $false, $true, 1, $false, "106×7=742", 2) | Out-Null
$d.Content.Find.Execute("651×9=5859", $true, $false, $false, $false, $false, $true, 1, $false, "548×8=4384", 2) | Out-Null
$d.Content.Find.Execute("301×4=1204", $true, $false, $false, $false, $false, $true, 1, $false, "826×4=3304", 2) | Out-Null
$d.Content.Find.Execute("307×6=1842", $true, $false, $false, $false, $false, $true, 1, $false, "526×7=3682", 2) | Out-Null
$d.Content.Find.Execute("658×4=2632", $true, $false, $false, $false, $false, $true, 1, $false, "832×8=6656", 2) | Out-Null
$d.Content.Find.Execute("692×7=4844", $true, $false, $false, $false, $false, $true, 1, $false, "844×3=2532", 2) | Out-Null
$d.Content.Find.Execute("451×7=3157", $true, $false, $false, $false, $false, $true, 1, $false, "607×6=3642", 2) | Out-Null
$d.Content.Find.Execute("619×5=3095", $true, $false, $false, $false, $false, $true, 1, $false, "639×4=2556", 2) | Out-Null
$d.Content.Find.Execute("187×6=1122", $true, $false, $false, $false, $false, $true, 1, $false, "880×6=5280", 2) | Out-Null
$d.Content.Find.Execute("418×9=3762", $true, $false, $false, $false, $false, $true, 1, $false, "865×3=2595", 2) | Out-Null
$d.Content.Find.Execute("237×4=948", $true, $false, $false, $false, $false, $true, 1, $false, "597×7=4179", 2) | Out-Null
$d.Content.Find.Execute("940×6=5640", $true, $false, $false, $false, $false, $true, 1, $false, "183×6=1098", 2) | Out-Null
$d.Content.Find.Execute("968×6=5808", $true, $false, $false, $false, $false, $true, 1, $false, "494×9=4446", 2) | Out-Null
$d.Content.Find.Execute("140×7=980", $true, $false, $false, $false, $false, $true, 1, $false, "455×6=2730", 2) | Out-Null
